$p = $ppt.ActivePresentation

# 1) Update the "datetimeFigureOut" date placeholder text on the Slide Master
#    and on every slide layout from "2011-12-04" to "2011/12/27".
$sm = $p.SlideMaster

for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $shp = $sm.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "2011-12-04") {
            $shp.TextFrame.TextRange.Text = "2011/12/27"
        }
    }
}

for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $cl = $sm.CustomLayouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $shp = $cl.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "2011-12-04") {
                $shp.TextFrame.TextRange.Text = "2011/12/27"
            }
        }
    }
}

# 2) Rename the "Physical Data Model" rounded-rectangle label to
#    "Template Data Object" on slide 1 (nested inside the group shape).
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(2)
$target = $grp.GroupItems.Item(6)
if ($target.TextFrame.TextRange.Text -eq "Physical Data Model") {
    $target.TextFrame.TextRange.Text = "Template Data Object"
}
